$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CVR column (A) holds digit-only values but must stay text (shared string),
# matching existing rows 2-6. Temporarily mark as Text format so COM doesn't
# auto-coerce to a number, then restore the default "Normal" style.
$ws.Range("A7:A9").NumberFormat = "@"

# ---- Row 7: new 2024Q1 cancellation ----
$ws.Range("A7").Value = "20246693"
$ws.Range("B7").Value = 2024
$ws.Range("C7").Value = 186200
$ws.Range("D7").Value = "Visma Løn"
$ws.Range("E7").Value = 45373
$ws.Range("E7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G7").Value = "Lessor"
$ws.Range("H7").Value = "2024Q1"
$ws.Range("I7").Value = "180000-200000"

# ---- Row 8: new 2024Q3 cancellation ----
$ws.Range("A8").Value = "15223332"
$ws.Range("B8").Value = 2024
$ws.Range("C8").Value = 191400
$ws.Range("D8").Value = "Visma Løn og HR"
$ws.Range("E8").Value = 45476
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H8").Value = "2024Q3"
$ws.Range("I8").Value = "180000-200000"

# ---- Row 9: new 2024Q3 cancellation ----
$ws.Range("A9").Value = "43268570"
$ws.Range("B9").Value = 2024
$ws.Range("C9").Value = 193860
$ws.Range("D9").Value = "BPO Løn & HR"
$ws.Range("E9").Value = 45526
$ws.Range("E9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H9").Value = "2024Q3"
$ws.Range("I9").Value = "180000-200000"

# Restore the CVR cells to the default (unstyled) look, like A2:A6.
$ws.Range("A7:A9").Style = "Normal"
